$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the roman-numeral "I" to the digit "1" in the unit name "Iº RCGd" -> "1º RCGd"
$ws.Range("A17").Value = "1º RCGd"

# Update the active selection to the cell that was last edited
$ws.Range("E7").Select()
